$p = $ppt.ActivePresentation

# The new "PART IV" divider slide is a copy of the existing "PART III"
# divider slide (slide 4), moved to the very end of the deck.
$src = $p.Slides.Item(4)
$new = $src.Duplicate()
$new.Item(1).MoveTo($p.Slides.Count)

$slide = $p.Slides.Item($p.Slides.Count)
$slide.Shapes.Item(1).TextFrame.TextRange.Text = "PART IV"
